$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume number and date range) ---
$ws.Range("A8").Characters(21,2).Text = "30"
$ws.Range("C9").Characters(27,9).Text = "7/21/2025"
$ws.Range("C9").Characters(47,9).Text = "7/27/2025"

# --- Cells changing data type / style (copy template cell then set value) ---
# Template cells (unchanged elsewhere in the diff):
#   D15 = style 13, text "0"
#   E15 = style 13, text "***.*"
#   G15 = style 14, plain number
#   K15 = style 15, plain decimal number
$ws.Range("D15").Copy($ws.Range("F15"))
$ws.Range("G15").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 3
$ws.Range("K15").Copy($ws.Range("E17"))
$ws.Range("E17").Value = 0
$ws.Range("D15").Copy($ws.Range("C20"))
$ws.Range("D15").Copy($ws.Range("C22"))
$ws.Range("D15").Copy($ws.Range("D22"))
$ws.Range("E15").Copy($ws.Range("E22"))
$ws.Range("D15").Copy($ws.Range("C23"))
$ws.Range("D15").Copy($ws.Range("C27"))

# --- Plain numeric value updates ---
$ws.Range("H15").Value = -100
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 140
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 60
$ws.Range("L16").Value = -25
$ws.Range("M16").Value = 9.090909090909
$ws.Range("N16").Value = -83.651226158038
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 88
$ws.Range("J17").Value = 68
$ws.Range("K17").Value = 29.411764705882
$ws.Range("L17").Value = -13.725490196078
$ws.Range("M17").Value = 6.024096385542
$ws.Range("N17").Value = -39.310344827586
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 59
$ws.Range("J18").Value = 61
$ws.Range("K18").Value = -3.278688524590
$ws.Range("L18").Value = -14.492753623188
$ws.Range("M18").Value = 1.724137931034
$ws.Range("N18").Value = -81.504702194357
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -68.75
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = -34.426229508196
$ws.Range("I19").Value = 341
$ws.Range("J19").Value = 388
$ws.Range("K19").Value = -12.113402061855
$ws.Range("L19").Value = -23.198198198198
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -25.869565217391
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 0
$ws.Range("L20").Value = -69.444444444444
$ws.Range("M20").Value = -56
$ws.Range("N20").Value = -96
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -41.666666666666
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = -19.318181818181
$ws.Range("I21").Value = 571
$ws.Range("J21").Value = 612
$ws.Range("K21").Value = -6.699346405228
$ws.Range("L21").Value = -22.207084468664
$ws.Range("M21").Value = 0.351493848857
$ws.Range("N21").Value = -63.974763406940
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = -3.846153846153
$ws.Range("L23").Value = -10.714285714285
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -27.777777777777
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 50
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 477
$ws.Range("J24").Value = 377
$ws.Range("K24").Value = 26.525198938992
$ws.Range("L24").Value = 8.656036446469
$ws.Range("M24").Value = -4.790419161676
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 64.285714285714
$ws.Range("I25").Value = 258
$ws.Range("J25").Value = 156
$ws.Range("K25").Value = 65.384615384615
$ws.Range("L25").Value = -1.526717557251
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 35
$ws.Range("H26").Value = 66.666666666666
$ws.Range("I26").Value = 173
$ws.Range("J26").Value = 182
$ws.Range("K26").Value = -4.945054945054
$ws.Range("L26").Value = -13.5
$ws.Range("M26").Value = -12.182741116751
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("C28").Value = 4
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 40
$ws.Range("I28").Value = 32
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = 3.225806451612
$ws.Range("L28").Value = -11.111111111111
